# Update the "Förändrad" (Changed) date column (C) for rows 2-89
# from serial date 45184 to 45185 (i.e. +1 day), matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 89
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45185
    }
}
